$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 282.7879796666667
$ws.Range("H2").Value = 848.363939
$ws.Range("I2").Value = 0.9674521741401267
$ws.Range("J2").Value = 0.9674521741401266
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.126464333333334
$ws.Range("N2").Value = 18.379393
$ws.Range("O2").Value = 0.1081098818071741
$ws.Range("P2").Value = 0.1081098818071741
$ws.Range("Q2").Value = 1732.490471323225
$ws.Range("R2").Value = 15592.41424190903
$ws.Range("S2").Value = 0.1045911402003827
$ws.Range("T2").Value = 0.1045911402003827

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 282.7879796666667
$ws.Range("H3").Value = 848.363939
$ws.Range("I3").Value = 0.9674521741401267
$ws.Range("J3").Value = 0.9674521741401266
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.73736866666667
$ws.Range("N3").Value = 107.212106
$ws.Range("O3").Value = 0.6306349784216607
$ws.Range("P3").Value = 0.6306349784216608
$ws.Range("Q3").Value = 10106.0982838495
$ws.Range("R3").Value = 90954.88455464554
$ws.Range("S3").Value = 0.6101091809628476
$ws.Range("T3").Value = 0.6101091809628476

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 282.7879796666667
$ws.Range("H4").Value = 848.363939
$ws.Range("I4").Value = 0.9674521741401267
$ws.Range("J4").Value = 0.9674521741401266
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.80503233333333
$ws.Range("N4").Value = 44.415097
$ws.Range("O4").Value = 0.2612551397711651
$ws.Range("P4").Value = 0.2612551397711651
$ws.Range("Q4").Value = 4186.68518244301
$ws.Range("R4").Value = 37680.16664198708
$ws.Range("S4").Value = 0.2527518529768964
$ws.Range("T4").Value = 0.2527518529768963

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.714696666666668
$ws.Range("H5").Value = 23.14409
$ws.Range("I5").Value = 0.02639291836872237
$ws.Range("J5").Value = 0.02639291836872237
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.126464333333334
$ws.Range("N5").Value = 18.379393
$ws.Range("O5").Value = 0.1081098818071741
$ws.Range("P5").Value = 0.1081098818071741
$ws.Range("Q5").Value = 47.26381397081889
$ws.Range("R5").Value = 425.37432573737
$ws.Range("S5").Value = 0.002853335285388971
$ws.Range("T5").Value = 0.002853335285388971

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.714696666666668
$ws.Range("H6").Value = 23.14409
$ws.Range("I6").Value = 0.02639291836872237
$ws.Range("J6").Value = 0.02639291836872237
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 35.73736866666667
$ws.Range("N6").Value = 107.212106
$ws.Range("O6").Value = 0.6306349784216607
$ws.Range("P6").Value = 0.6306349784216608
$ws.Range("Q6").Value = 275.7029589281711
$ws.Range("R6").Value = 2481.32663035354
$ws.Range("S6").Value = 0.01664429750594389
$ws.Range("T6").Value = 0.01664429750594389

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.714696666666668
$ws.Range("H7").Value = 23.14409
$ws.Range("I7").Value = 0.02639291836872237
$ws.Range("J7").Value = 0.02639291836872237
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.80503233333333
$ws.Range("N7").Value = 44.415097
$ws.Range("O7").Value = 0.2612551397711651
$ws.Range("P7").Value = 0.2612551397711651
$ws.Range("Q7").Value = 114.2163335918589
$ws.Range("R7").Value = 1027.94700232673
$ws.Range("S7").Value = 0.006895285577389515
$ws.Range("T7").Value = 0.006895285577389514

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.799090333333333
$ws.Range("H8").Value = 5.397271
$ws.Range("I8").Value = 0.006154907491150983
$ws.Range("J8").Value = 0.006154907491150983
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.126464333333334
$ws.Range("N8").Value = 18.379393
$ws.Range("O8").Value = 0.1081098818071741
$ws.Range("P8").Value = 0.1081098818071741
$ws.Range("Q8").Value = 11.02206275961144
$ws.Range("R8").Value = 99.198564836503
$ws.Range("S8").Value = 0.0006654063214024233
$ws.Range("T8").Value = 0.0006654063214024233

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.799090333333333
$ws.Range("H9").Value = 5.397271
$ws.Range("I9").Value = 0.006154907491150983
$ws.Range("J9").Value = 0.006154907491150983
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 35.73736866666667
$ws.Range("N9").Value = 107.212106
$ws.Range("O9").Value = 0.6306349784216607
$ws.Range("P9").Value = 0.6306349784216608
$ws.Range("Q9").Value = 64.29475450696955
$ws.Range("R9").Value = 578.6527905627261
$ws.Range("S9").Value = 0.003881499952869318
$ws.Range("T9").Value = 0.003881499952869319

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.799090333333333
$ws.Range("H10").Value = 5.397271
$ws.Range("I10").Value = 0.006154907491150983
$ws.Range("J10").Value = 0.006154907491150983
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.80503233333333
$ws.Range("N10").Value = 44.415097
$ws.Range("O10").Value = 0.2612551397711651
$ws.Range("P10").Value = 0.2612551397711651
$ws.Range("Q10").Value = 26.63559055558745
$ws.Range("R10").Value = 239.720315000287
$ws.Range("S10").Value = 0.001608001216879241
$ws.Range("T10").Value = 0.001608001216879241

